$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "328.95"
Set-TextValue "E2" "6.61%"
Set-TextValue "D3" "40.60"
Set-TextValue "E3" "11.70%"
Set-TextValue "D4" "5.941"
Set-TextValue "E4" "16.09%"
Set-TextValue "D5" "0.08153"
Set-TextValue "E5" "5.62%"
Set-TextValue "D6" "4.565"
Set-TextValue "E6" "4.11%"
Set-TextValue "D7" "8.743"
Set-TextValue "E7" "5.42%"
Set-TextValue "D8" "1.948"
Set-TextValue "E8" "4.92%"
Set-TextValue "E9" "-1.17%"
Set-TextValue "D10" "0.9455"
Set-TextValue "E10" "2.66%"
Set-TextValue "D11" "0.1310"
Set-TextValue "E11" "15.69%"
Set-TextValue "D12" "0.1996"
Set-TextValue "E12" "7.58%"
Set-TextValue "D13" "0.09292"
Set-TextValue "E13" "5.92%"
Set-TextValue "D14" "0.03432"
Set-TextValue "E14" "3.25%"
Set-TextValue "D15" "0.09626"
Set-TextValue "E15" "1.02%"
Set-TextValue "D16" "0.001331"
Set-TextValue "E16" "-3.38%"
Set-TextValue "D17" "0.006019"
Set-TextValue "E17" "-2.04%"
Set-TextValue "D18" "3.375"
Set-TextValue "E18" "0.32%"
Set-TextValue "D19" "0.3495"
Set-TextValue "E19" "1.45%"
Set-TextValue "D20" "7.740"
Set-TextValue "E20" "22.55%"
Set-TextValue "D21" "0.1438"
Set-TextValue "D22" "0.2450"
Set-TextValue "E22" "5.87%"
Set-TextValue "D23" "0.04428"
Set-TextValue "E23" "2.13%"
Set-TextValue "E24" "4.21%"
Set-TextValue "D25" "0.004358"
Set-TextValue "E25" "2.28%"
Set-TextValue "D26" "0.0001191"
Set-TextValue "E26" "-10.62%"
Set-TextValue "D27" "0.0003991"
Set-TextValue "E27" "37.41%"
Set-TextValue "D39" "0.02504"
Set-TextValue "E39" "19.13%"
Set-TextValue "D40" "0.05305"
Set-TextValue "E40" "7.68%"
Set-TextValue "D41" "0.007606"
Set-TextValue "E41" "0.52%"
Set-TextValue "D42" "0.1433"
Set-TextValue "E42" "6.24%"
Set-TextValue "D43" "0.008923"
Set-TextValue "E43" "4.17%"
Set-TextValue "E44" "-0.38%"
Set-TextValue "D45" "0.009500"
Set-TextValue "D46" "0.00006864"
Set-TextValue "E46" "6.12%"
Set-TextValue "E47" "-0.11%"
Set-TextValue "D48" "0.002899"
Set-TextValue "E48" "-12.17%"
Set-TextValue "E49" "24.60%"
Set-TextValue "E50" "-0.11%"
Set-TextValue "E51" "-0.11%"
